# Auto-generated Excel COM-interop script
# Updates Universalis market-price-derived columns (H-N) across the 8 Anima Profits sheets
# to reflect a scheduled market-data refresh, per the commit's unified diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 880.6338
$ws.Range("I17").Value = 561.96155
$ws.Range("J17").Value = 1064.7556
$ws.Range("K17").Value = 1685.88465
$ws.Range("L17").Value = 3194.2668
$ws.Range("M17").Value = -1517.88465
$ws.Range("N17").Value = -3530.2668

$ws.Range("H40").Value = 1457.4286

$ws.Range("H112").Value = 5178.6113
$ws.Range("J112").Value = 5616.0605
$ws.Range("L112").Value = 16848.1815
$ws.Range("N112").Value = -19064.1815

$ws.Range("H129").Value = 1320.9333
$ws.Range("J129").Value = 1665.7778
$ws.Range("L129").Value = 4997.3334
$ws.Range("N129").Value = -14997.3334

$ws.Range("H132").Value = 2602.3262
$ws.Range("I132").Value = 2430.4187
$ws.Range("K132").Value = 7291.256100000001
$ws.Range("M132").Value = -4761.256100000001

$ws.Range("H135").Value = 1301.8529
$ws.Range("I135").Value = 508.21875
$ws.Range("J135").Value = 14000
$ws.Range("K135").Value = 4573.96875
$ws.Range("L135").Value = 126000
$ws.Range("M135").Value = -2038.96875
$ws.Range("N135").Value = -131070

$ws.Range("H141").Value = 4591.7856
$ws.Range("I141").Value = 1610.5883
$ws.Range("J141").Value = 9199.091
$ws.Range("K141").Value = 4831.7649
$ws.Range("L141").Value = 27597.273
$ws.Range("M141").Value = 348.2350999999999
$ws.Range("N141").Value = -37957.273


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 156.16667
$ws.Range("I5").Value = 147.4
$ws.Range("K5").Value = 147.4
$ws.Range("M5").Value = -35.40000000000001

$ws.Range("H32").Value = 835434.75
$ws.Range("I32").Value = 896190.2
$ws.Range("J32").Value = 33462.8
$ws.Range("K32").Value = 896190.2
$ws.Range("L32").Value = 33462.8
$ws.Range("M32").Value = -895903.2
$ws.Range("N32").Value = -34036.8

$ws.Range("H98").Value = 24000
$ws.Range("J98").Value = 24000
$ws.Range("L98").Value = 24000
$ws.Range("N98").Value = -29990


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 156.16667
$ws.Range("I4").Value = 147.4
$ws.Range("K4").Value = 147.4
$ws.Range("M4").Value = -32.40000000000001

$ws.Range("H22").Value = 279.5238
$ws.Range("I22").Value = 279.5238
$ws.Range("K22").Value = 279.5238
$ws.Range("M22").Value = -106.5238

$ws.Range("H81").Value = 40519
$ws.Range("J81").Value = 40519
$ws.Range("L81").Value = 40519
$ws.Range("N81").Value = -42641

$ws.Range("H84").Value = 40519
$ws.Range("J84").Value = 40519
$ws.Range("L84").Value = 121557
$ws.Range("N84").Value = -132165


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 51666.668
$ws.Range("J28").Value = 51666.668
$ws.Range("L28").Value = 51666.668
$ws.Range("N28").Value = -52156.668

$ws.Range("H31").Value = 3192.7285
$ws.Range("I31").Value = 943.88464
$ws.Range("J31").Value = 7225.1377
$ws.Range("K31").Value = 943.88464
$ws.Range("L31").Value = 7225.1377
$ws.Range("M31").Value = -648.88464
$ws.Range("N31").Value = -7815.1377

$ws.Range("H34").Value = 3192.7285
$ws.Range("I34").Value = 943.88464
$ws.Range("J34").Value = 7225.1377
$ws.Range("K34").Value = 943.88464
$ws.Range("L34").Value = 7225.1377
$ws.Range("M34").Value = -741.88464
$ws.Range("N34").Value = -7629.1377

$ws.Range("H43").Value = 117000
$ws.Range("J43").Value = 117000
$ws.Range("L43").Value = 117000
$ws.Range("N43").Value = -117368

$ws.Range("H81").Value = 54994.668
$ws.Range("J81").Value = 54994.668
$ws.Range("L81").Value = 54994.668
$ws.Range("N81").Value = -56990.668

$ws.Range("H84").Value = 54994.668
$ws.Range("J84").Value = 54994.668
$ws.Range("L84").Value = 164984.004
$ws.Range("N84").Value = -174968.004

$ws.Range("H86").Value = 2094
$ws.Range("I86").Value = 2257.4443
$ws.Range("J86").Value = 1726.25
$ws.Range("K86").Value = 2257.4443
$ws.Range("L86").Value = 1726.25
$ws.Range("M86").Value = -1134.4443
$ws.Range("N86").Value = -3972.25

$ws.Range("H89").Value = 2094
$ws.Range("I89").Value = 2257.4443
$ws.Range("J89").Value = 1726.25
$ws.Range("K89").Value = 11287.2215
$ws.Range("L89").Value = 8631.25
$ws.Range("M89").Value = -5671.2215
$ws.Range("N89").Value = -19863.25

$ws.Range("H101").Value = 117000
$ws.Range("J101").Value = 117000
$ws.Range("L101").Value = 117000
$ws.Range("N101").Value = -123490

$ws.Range("H132").Value = 1931.8611
$ws.Range("I132").Value = 1768.1111
$ws.Range("K132").Value = 5304.3333
$ws.Range("M132").Value = -2774.3333

$ws.Range("H140").Value = 73640
$ws.Range("J140").Value = 74853.336
$ws.Range("L140").Value = 74853.336
$ws.Range("N140").Value = -85213.336


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 117.77778
$ws.Range("J12").Value = 145.71428
$ws.Range("L12").Value = 437.14284
$ws.Range("N12").Value = -783.14284

$ws.Range("H33").Value = 71614.5
$ws.Range("J33").Value = 450050
$ws.Range("L33").Value = 2700300
$ws.Range("N33").Value = -2700866

$ws.Range("H68").Value = 8281
$ws.Range("I68").Value = 45900
$ws.Range("J68").Value = 1441.1818
$ws.Range("K68").Value = 137700
$ws.Range("L68").Value = 4323.5454
$ws.Range("M68").Value = -136889
$ws.Range("N68").Value = -5945.5454

$ws.Range("H71").Value = 8281
$ws.Range("I71").Value = 45900
$ws.Range("J71").Value = 1441.1818
$ws.Range("K71").Value = 413100
$ws.Range("L71").Value = 12970.6362
$ws.Range("M71").Value = -409044
$ws.Range("N71").Value = -21082.6362

$ws.Range("H74").Value = 2833.3333
$ws.Range("J74").Value = 3750
$ws.Range("L74").Value = 11250
$ws.Range("N74").Value = -13372

$ws.Range("H77").Value = 2833.3333
$ws.Range("J77").Value = 3750
$ws.Range("L77").Value = 33750
$ws.Range("N77").Value = -44358

$ws.Range("H113").Value = 546.4231
$ws.Range("I113").Value = 568.2
$ws.Range("J113").Value = 516.7273
$ws.Range("K113").Value = 1704.6
$ws.Range("L113").Value = 1550.1819
$ws.Range("M113").Value = 465.3999999999999
$ws.Range("N113").Value = -5890.1819

$ws.Range("H120").Value = 10257.5
$ws.Range("J120").Value = 13000
$ws.Range("L120").Value = 39000
$ws.Range("N120").Value = -48676

$ws.Range("H122").Value = 3977.9355
$ws.Range("J122").Value = 8823.154
$ws.Range("L122").Value = 79408.386
$ws.Range("N122").Value = -84308.386

$ws.Range("H131").Value = 2786.3484
$ws.Range("I131").Value = 469.3
$ws.Range("J131").Value = 3200.1072
$ws.Range("K131").Value = 1407.9
$ws.Range("L131").Value = 9600.321599999999
$ws.Range("M131").Value = 3632.1
$ws.Range("N131").Value = -19680.3216

$ws.Range("H136").Value = 3859.95
$ws.Range("I136").Value = 1333.1666
$ws.Range("J136").Value = 4942.857
$ws.Range("K136").Value = 3999.4998
$ws.Range("L136").Value = 14828.571
$ws.Range("M136").Value = 1100.5002
$ws.Range("N136").Value = -25028.571

$ws.Range("H137").Value = 6951882
$ws.Range("J137").Value = 3922.2222
$ws.Range("L137").Value = 11766.6666
$ws.Range("N137").Value = -21966.6666

$ws.Range("H139").Value = 2588.585
$ws.Range("J139").Value = 3447.6128
$ws.Range("L139").Value = 10342.8384
$ws.Range("N139").Value = -20622.8384

$ws.Range("H140").Value = 1563.2354
$ws.Range("I140").Value = 1129.5454
$ws.Range("K140").Value = 3388.6362
$ws.Range("M140").Value = 1791.3638


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 335334.66
$ws.Range("I40").Value = 501002
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 501002
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -500866
$ws.Range("N40").Value = -4272

$ws.Range("H132").Value = 2362.2456
$ws.Range("I132").Value = 2090.1353
$ws.Range("J132").Value = 2865.65
$ws.Range("K132").Value = 6270.4059
$ws.Range("L132").Value = 8596.950000000001
$ws.Range("M132").Value = -3740.4059
$ws.Range("N132").Value = -13656.95

$ws.Range("H136").Value = 11113965
$ws.Range("I136").Value = 3834
$ws.Range("K136").Value = 11502
$ws.Range("M136").Value = -8952


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 70004
$ws.Range("J3").Value = 70004
$ws.Range("L3").Value = 70004
$ws.Range("N3").Value = -70232

$ws.Range("H124").Value = 54107.25
$ws.Range("J124").Value = 54107.25
$ws.Range("L124").Value = 54107.25
$ws.Range("N124").Value = -63927.25

$ws.Range("H136").Value = 3186.5454
$ws.Range("I136").Value = 2790.4707
$ws.Range("K136").Value = 8371.4121
$ws.Range("M136").Value = -5821.4121

